$wb = $excel.ActiveWorkbook

# Sheet1 ("PadaliHMSfinal" small table): selection now covers the whole used range
$sheet1 = $wb.Worksheets.Item("Sheet1")
[void]$sheet1.Range("A1:G13").Select()

# Add the new "Sheet2" worksheet (appended after current sheets) with the
# Kanhur mesai village monthly weather-budget data.
$ws = $wb.Worksheets.Add()
$ws.Name = "Sheet2"

$ws.Cells.Item(1,1).Value = "Month"
$ws.Cells.Item(1,2).Value = "Evaporation"
$ws.Cells.Item(1,3).Value = "Humidity"
$ws.Cells.Item(1,4).Value = "Rainfall"
$ws.Cells.Item(1,5).Value = "Min tem"
$ws.Cells.Item(1,6).Value = "Max tem"
$ws.Cells.Item(1,7).Value = "Wind speed"

$ws.Cells.Item(2,1).Value = "January"
$ws.Cells.Item(2,2).Value = 0
$ws.Cells.Item(2,3).Value = 0
$ws.Cells.Item(2,4).Value = 0
$ws.Cells.Item(2,5).Value = 18
$ws.Cells.Item(2,6).Value = 29
$ws.Cells.Item(2,7).Value = 8.7

$ws.Cells.Item(3,1).Value = "February"
$ws.Cells.Item(3,2).Value = 0
$ws.Cells.Item(3,3).Value = 0
$ws.Cells.Item(3,4).Value = 0
$ws.Cells.Item(3,5).Value = 20
$ws.Cells.Item(3,6).Value = 33
$ws.Cells.Item(3,7).Value = 8.9

$ws.Cells.Item(4,1).Value = "March"
$ws.Cells.Item(4,2).Value = 0
$ws.Cells.Item(4,3).Value = 0
$ws.Cells.Item(4,4).Value = 0
$ws.Cells.Item(4,5).Value = 20
$ws.Cells.Item(4,6).Value = 34
$ws.Cells.Item(4,7).Value = 11.1

$ws.Cells.Item(5,1).Value = "April"
$ws.Cells.Item(5,2).Value = 0
$ws.Cells.Item(5,3).Value = 0
$ws.Cells.Item(5,4).Value = 0
$ws.Cells.Item(5,5).Value = 24
$ws.Cells.Item(5,6).Value = 39
$ws.Cells.Item(5,7).Value = 11.3

$ws.Cells.Item(6,1).Value = "May"
$ws.Cells.Item(6,2).Value = 0
$ws.Cells.Item(6,3).Value = 0
$ws.Cells.Item(6,4).Value = 0
$ws.Cells.Item(6,5).Value = 23
$ws.Cells.Item(6,6).Value = 39
$ws.Cells.Item(6,7).Value = 14.4

$ws.Cells.Item(7,1).Value = "June"
$ws.Cells.Item(7,2).Value = 0
$ws.Cells.Item(7,3).Value = 0
$ws.Cells.Item(7,4).Value = 0
$ws.Cells.Item(7,5).Value = 22
$ws.Cells.Item(7,6).Value = 29
$ws.Cells.Item(7,7).Value = 14.9

$ws.Cells.Item(8,1).Value = "July"
$ws.Cells.Item(8,2).Value = 0
$ws.Cells.Item(8,3).Value = 0
$ws.Cells.Item(8,4).Value = 0
$ws.Cells.Item(8,5).Value = 22
$ws.Cells.Item(8,6).Value = 28
$ws.Cells.Item(8,7).Value = 15.2

$ws.Cells.Item(9,1).Value = "August"
$ws.Cells.Item(9,2).Value = 0
$ws.Cells.Item(9,3).Value = 0
$ws.Cells.Item(9,4).Value = 0
$ws.Cells.Item(9,5).Value = 22
$ws.Cells.Item(9,6).Value = 26
$ws.Cells.Item(9,7).Value = 18.3

$ws.Cells.Item(10,1).Value = "September"
$ws.Cells.Item(10,2).Value = 0
$ws.Cells.Item(10,3).Value = 0
$ws.Cells.Item(10,4).Value = 0
$ws.Cells.Item(10,5).Value = 22
$ws.Cells.Item(10,6).Value = 28
$ws.Cells.Item(10,7).Value = 9.6

$ws.Cells.Item(11,1).Value = "October"
$ws.Cells.Item(11,2).Value = 0
$ws.Cells.Item(11,3).Value = 0
$ws.Cells.Item(11,4).Value = 0
$ws.Cells.Item(11,5).Value = 22
$ws.Cells.Item(11,6).Value = 29
$ws.Cells.Item(11,7).Value = 7

$ws.Cells.Item(12,1).Value = "November"
$ws.Cells.Item(12,2).Value = 0
$ws.Cells.Item(12,3).Value = 0
$ws.Cells.Item(12,4).Value = 0
$ws.Cells.Item(12,5).Value = 20
$ws.Cells.Item(12,6).Value = 29
$ws.Cells.Item(12,7).Value = 9

$ws.Cells.Item(13,1).Value = "December"
$ws.Cells.Item(13,2).Value = 0
$ws.Cells.Item(13,3).Value = 0
$ws.Cells.Item(13,4).Value = 0
$ws.Cells.Item(13,5).Value = 20
$ws.Cells.Item(13,6).Value = 28
$ws.Cells.Item(13,7).Value = 8.3

# Move the new sheet to the end of the tab strip (after Sheet1 & Sheet3)
$ws.Move([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))

# Un-select Sheet3's tab (it was tabSelected before) and make the new
# Sheet2 ("Kanhur mesai" result shown to CARD) the active tab/selection.
[void]$wb.Worksheets.Item(3).Select()
[void]$wb.Worksheets.Item(3).Range("F13").Select()
